# 03项目计划表 - "Add files via upload" update
# Adds the completion data for the 2018.10.25 (Wed) week block (rows 123-128)
# and appends a brand-new week block (2018.10.25 Thu) as rows 131-140.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Fill in the "completion" column for the previous week block ---
$ws.Range("C123").Value = "完成"
$ws.Range("C124").Value = "完成"
$ws.Range("C125").Value = 1
$ws.Range("C126").Value = 1
$ws.Range("C127").Value = 1
$ws.Range("C128").Value = 0.8

# --- 2. Replace the placeholder summary text for that block ---
$ws.Range("A129").Value = "总结：我们现在交互还是这么不是很顺利，经常遇到困难，我们应该反省一下自己到底哪里错了，找一下自己的原因，这样才能走得更远啊，但是我们一定不能气馁，距离任务任重而道远，千万不要放弃。"

# --- 3. Build the new week block (rows 131-140) by cloning the formatting
#        of the previous block (rows 121-130), then filling in values ---
$ws.Range("A121:D130").Copy()
$ws.Range("A131").PasteSpecial(-4122)

$ws.Range("A131").Value = "日期：2018.10.25第九周周四"
$ws.Range("A131:D131").Merge()

$ws.Range("A133").Value = "陈柯赞"
$ws.Range("B133").Value = "pc端注册界面编码"

$ws.Range("A134").Value = "黎安生"
$ws.Range("B134").Value = "pc端登录界面编码"

$ws.Range("A135").Value = "王智永"
$ws.Range("B135").Value = "获取百度地图到底的位置信息"

$ws.Range("A136").Value = "郑海文"
$ws.Range("B136").Value = "申请加入群聊界面及请求"

$ws.Range("A137").Value = "赵华亮"
$ws.Range("B137").Value = "拒绝加入群聊界面及请求"

$ws.Range("A138").Value = "叶田"
$ws.Range("B138").Value = "找图片和素材"

$ws.Range("A139").Value = "总结："

$ws.Range("D133:D138").Merge()
$ws.Range("A139:D140").Merge()

$ws.Range("B136").Select()
